$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.482.96'
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").Value = '2.925.58'
$ws.Range("E3").Value = '  +4.38%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.73'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.16'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.40'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.16'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.90'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = '3.388.21'
$ws.Range("E15").Value = '  +4.59%  '
$ws.Range("D16").Value = '2.931.89'
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").Value = '52.511.30'
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.72'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  +5.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.50'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +6.71%  '
$ws.Range("D22").Value = '0.0₃0986'
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.23'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '271.83'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.81'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +2.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.14'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +4.13%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.67'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.09'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +2.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.54'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +5.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +9.17%  '
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '53.19'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0944'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +10.24%  '
$ws.Range("E36").Value = '  +2.07%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  +7.17%  '
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.91'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.09'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +4.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.75'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +10.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.59'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +13.18%  '
$ws.Range("E43").Value = '  +2.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.66'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("E45").Value = '  +7.60%  '
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.58'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +6.47%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.219.77'
$ws.Range("E48").Value = '  +4.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.265'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +25.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0343'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +16.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.964'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +6.39%  '
